$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows (by sheet row number) whose Grade column (B) should become "C"
$rows = @(11, 25, 36, 38, 59, 64, 67, 70, 74, 82, 93, 98)
foreach ($r in $rows) {
    $ws.Range("B$r").Value = "C"
}

# Update the view state: scrolled position and active selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F81").Select()
